# Regenerated localization-status report: two files moved from
# "Ready for handoff" back to "In Translation" (status regressed) for
# both the zh-cn and de-de locales. Update the per-locale status tables
# as well as the Overview roll-up sheet.
#
# Affected source files:
#   12c51228-8784-4454-bc7c-ae7d0be05400.md  (Overview row 3 / table row 2)
#   440d60bd-8fb4-4c83-9f59-c76cf27d4766.md  (Overview row 4 / table row 3)
#
# 796762e4-7f92-41be-85a6-414a4b2e9726.md stays "Ready for handoff".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: columns E (zh-cn) and F (de-de) hold the per-locale status
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus

# zh-cn sheet: column C is "Status"
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

# de-de sheet: column C is "Status"
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus
